$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 3416
$ws.Range("E3").Value = "24.4.2024"
$ws.Range("L3").Value = "Babichev Sergii, prof. CSc., DSc."
$ws.Range("B4").Value = 6973
$ws.Range("L4").Value = "Beránek Pavel, Ing. Mgr."
$ws.Range("B5").Value = 3416
$ws.Range("E5").Value = "15.5.2024"
$ws.Range("L5").Value = "Babichev Sergii, prof. CSc., DSc."
$ws.Range("B6").Value = 6259
$ws.Range("E6").Value = "17.4.2024"
$ws.Range("L6").Value = "Škvára Jiří, RNDr. Ph.D."
$ws.Range("B7").Value = 6973
$ws.Range("E7").Value = "17.4.2024"
$ws.Range("L7").Value = "Beránek Pavel, Ing. Mgr."
$ws.Range("B8").Value = 2220
$ws.Range("E8").Value = "24.4.2024"
$ws.Range("L8").Value = "Škvor Jiří, RNDr. Ph.D."
$ws.Range("E9").Value = "24.4.2024"
$ws.Range("E10").Value = "15.5.2024"
$ws.Range("B11").Value = 2220
$ws.Range("L11").Value = "Škvor Jiří, RNDr. Ph.D."
$ws.Range("B13").Value = 3416
$ws.Range("E13").Value = "17.4.2024"
$ws.Range("L13").Value = "Babichev Sergii, prof. CSc., DSc."
$ws.Range("B14").Value = 6973
$ws.Range("L14").Value = "Beránek Pavel, Ing. Mgr."
$ws.Range("E15").Value = "2.5.2024"
$ws.Range("E16").Value = "18.4.2024"
$ws.Range("E17").Value = "16.5.2024"
$ws.Range("E19").Value = "11.4.2024"
$ws.Range("E20").Value = "18.4.2024"
$ws.Range("E22").Value = "25.4.2024"
$ws.Range("E24").Value = "16.5.2024"
$ws.Range("E25").Value = "16.4.2024"
$ws.Range("E26").Value = "30.4.2024"
$ws.Range("E27").Value = "14.5.2024"
$ws.Range("E31").Value = "11.4.2024"
$ws.Range("E33").Value = "9.5.2024"
$ws.Range("E34").Value = "25.4.2024"
$ws.Range("E38").Value = "11.4.2024"
$ws.Range("E39").Value = "16.5.2024"
$ws.Range("E40").Value = "25.4.2024"
$ws.Range("E41").Value = "9.5.2024"
$ws.Range("E42").Value = "2.5.2024"
$ws.Range("E43").Value = "18.4.2024"
$ws.Range("E45").Value = "10.5.2024"
$ws.Range("E46").Value = "17.5.2024"
$ws.Range("E47").Value = "3.5.2024"
$ws.Range("E48").Value = "26.4.2024"
$ws.Range("E50").Value = "19.4.2024"
$ws.Range("E51").Value = "24.4.2024"
$ws.Range("E52").Value = "17.4.2024"
